$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell reference -> new text value, applied as plain text
# (preserving the original General number format / inlineStr semantics)
$updates = @{
    'D2' = '246.84'
    'E2' = '0.93%'
    'D3' = '29.55'
    'E3' = '8.40%'
    'D4' = '5.172'
    'E4' = '2.29%'
    'D5' = '0.05707'
    'D6' = '6.602'
    'E6' = '1.59%'
    'D7' = '0.8575'
    'E7' = '4.45%'
    'D8' = '0.8763'
    'E8' = '4.37%'
    'D9' = '0.1370'
    'E9' = '3.19%'
    'D10' = '0.07070'
    'E10' = '2.39%'
    'E11' = '0.56%'
    'D12' = '0.09383'
    'E12' = '-0.18%'
    'D13' = '0.001512'
    'E13' = '0.20%'
    'D14' = '0.04178'
    'E14' = '1.91%'
    'D15' = '0.0006017'
    'E15' = '-93.99%'
    'D16' = '0.006032'
    'E16' = '-3.26%'
    'E17' = '3,767.69%'
    'D18' = '3.485'
    'E18' = '-0.64%'
    'D19' = '3.080'
    'E19' = '2.67%'
    'D20' = '2.189'
    'E20' = '-5.34%'
    'E21' = '0.67%'
    'D22' = '0.03309'
    'E22' = '3.97%'
    'D23' = '0.1300'
    'E23' = '0.74%'
    'D24' = '3.470'
    'E24' = '-2.24%'
    'D25' = '0.1380'
    'E25' = '0.45%'
    'D26' = '0.005042'
    'E26' = '27.66%'
    'D27' = '0.001220'
    'E27' = '0.31%'
    'D28' = '0.0001210'
    'E28' = '23.52%'
    'E40' = '1.42%'
    'D41' = '0.005748'
    'E41' = '-3.32%'
    'E42' = '1.80%'
    'D43' = '0.002539'
    'E43' = '7.78%'
    'D44' = '0.01022'
    'E44' = '9.00%'
    'D45' = '0.00005157'
    'E45' = '-0.81%'
    'E46' = '0.02%'
    'D47' = '0.06997'
    'E47' = '-31.03%'
    'D48' = '0.002574'
    'E48' = '-0.28%'
    'E49' = '0.02%'
    'E50' = '0.02%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
